# Re-theme the deck: swap the "Integral"/Red Violet color scheme that is
# currently driving the slide master for the standard Office Theme color
# scheme (the palette that, before this edit, only lived on the Notes
# Master's theme part).
#
# The presentation's color scheme is reached through
# SlideMaster.Theme.ThemeColorScheme.Colors(n).RGB — each of the twelve
# theme colors (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink, in that
# COM index order) is pushed to the corresponding "Office" palette value.

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

function Set-ThemeColor {
    param($scheme, [int]$index, [int]$r, [int]$g, [int]$b)
    $scheme.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

# index : role      : target "Office" hex
# 1     : dk1        000000
# 2     : lt1        FFFFFF
# 3     : dk2        44546A
# 4     : lt2        E7E6E6
# 5     : accent1    5B9BD5
# 6     : accent2    ED7D31
# 7     : accent3    A5A5A5
# 8     : accent4    FFC000
# 9     : accent5    4472C4
# 10    : accent6    70AD47
# 11    : hlink      0563C1
# 12    : folHlink   954F72

Set-ThemeColor $tcs 1  0x00 0x00 0x00
Set-ThemeColor $tcs 2  0xFF 0xFF 0xFF
Set-ThemeColor $tcs 3  0x44 0x54 0x6A
Set-ThemeColor $tcs 4  0xE7 0xE6 0xE6
Set-ThemeColor $tcs 5  0x5B 0x9B 0xD5
Set-ThemeColor $tcs 6  0xED 0x7D 0x31
Set-ThemeColor $tcs 7  0xA5 0xA5 0xA5
Set-ThemeColor $tcs 8  0xFF 0xC0 0x00
Set-ThemeColor $tcs 9  0x44 0x72 0xC4
Set-ThemeColor $tcs 10 0x70 0xAD 0x47
Set-ThemeColor $tcs 11 0x05 0x63 0xC1
Set-ThemeColor $tcs 12 0x95 0x4F 0x72

Write-Output "Theme colors updated to Office palette"
